$wb = $excel.ActiveWorkbook

# --- Hoja "1er Parcial" ---
$ws1 = $wb.Worksheets.Item("1er Parcial")

$ws1.Range("E2").Value = 19
$ws1.Range("F2").Value = 7
$ws1.Range("G2").Value = 73.08
$ws1.Range("H2").Value = 26.92
$ws1.Range("I2").Value = 7.7
$ws1.Range("J2").Value = 7
$ws1.Range("K2").Value = 26.92

$ws1.Range("E3").Value = 20
$ws1.Range("F3").Value = 6
$ws1.Range("G3").Value = 76.92
$ws1.Range("H3").Value = 23.08
$ws1.Range("I3").Value = 7.7
$ws1.Range("J3").Value = 6
$ws1.Range("K3").Value = 23.08

# --- Hoja "2o Parcial" ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

$ws2.Range("E2").Value = 14
$ws2.Range("F2").Value = 12
$ws2.Range("G2").Value = 53.85
$ws2.Range("H2").Value = 46.15
$ws2.Range("I2").Value = 8.800000000000001
$ws2.Range("J2").Value = 12
$ws2.Range("K2").Value = 46.15

$ws2.Range("E3").Value = 14
$ws2.Range("F3").Value = 12
$ws2.Range("G3").Value = 53.85
$ws2.Range("H3").Value = 46.15
$ws2.Range("I3").Value = 8.199999999999999
$ws2.Range("J3").Value = 12
$ws2.Range("K3").Value = 46.15

# --- Hoja "3er Parcial" ---
$ws3 = $wb.Worksheets.Item("3er Parcial")

$ws3.Range("E2").Value = 19
$ws3.Range("F2").Value = 7
$ws3.Range("G2").Value = 73.08
$ws3.Range("H2").Value = 26.92
$ws3.Range("I2").Value = 8.1
$ws3.Range("J2").Value = 7
$ws3.Range("K2").Value = 26.92

$ws3.Range("E3").Value = 20
$ws3.Range("F3").Value = 6
$ws3.Range("G3").Value = 76.92
$ws3.Range("H3").Value = 23.08
$ws3.Range("J3").Value = 6
$ws3.Range("K3").Value = 23.08
